# README for basic usage
# Rewrites the "combo" sheet's sample data (adds a 3rd [Type]/[Amount] pair,
# renames the header keys from the old a.b.c.d.e.f / a.d.c.d.e.f. scheme to
# the new a.b / a.b3. scheme, and swaps the sample row's literal values),
# and moves the saved cursor/selection on a few sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # simple
$ws2 = $wb.Worksheets.Item(2)   # vector
$ws3 = $wb.Worksheets.Item(3)   # vector_dict
$ws4 = $wb.Worksheets.Item(4)   # combo

# --- "combo" sheet data rewrite -------------------------------------------
# Populate the new values in the same order the original author touched the
# cells in (some row 2 string cells first, then the header row, then the
# rest of row 2) so shared-string interning lines up with the saved file.
$ws4.Range("C2").Value = "Cake"
$ws4.Range("E2").Value = "Chocolate"
$ws4.Range("H2").Value = "Salt;100"
$ws4.Range("G2").Value = "Ingredient"

$ws4.Range("A1").Value = "a.b"
$ws4.Range("B1").Value = "a.b2(Int)"
$ws4.Range("C1").Value = "a.b3.[1,Type]"
$ws4.Range("D1").Value = "a.b3.[1,Amount]"
$ws4.Range("E1").Value = "a.b3.[2,Type]"
$ws4.Range("F1").Value = "a.b3.[2,Amount]"
$ws4.Range("G1").Value = "a.b3.[3,Type]"
$ws4.Range("H1").Value = "a.b3.[3,Amount()]"

$ws4.Range("A2").Value = "Fooood"
$ws4.Range("B2").Value = "100;200;300"
$ws4.Range("D2").Value = 50
$ws4.Range("F2").Value = 19

# New column H needs the same cell style the rest of the table uses.
$ws4.Range("G1:G2").Copy() | Out-Null
$ws4.Range("H1:H2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- saved cursor / selection per sheet ------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("E1").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("F1").Select() | Out-Null

$ws4.Activate() | Out-Null
$ws4.Range("A3").Select() | Out-Null

# --- page setup on the "combo" sheet ---------------------------------------
$ps = $ws4.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
